$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.775549333333333
$ws.Range("H2").Value = 20.326648
$ws.Range("I2").Value = 0.5307754563424079
$ws.Range("J2").Value = 0.5307754563424079
$ws.Range("M2").Value = 70.23436
$ws.Range("N2").Value = 210.70308
$ws.Range("O2").Value = 0.7023186840741513
$ws.Range("P2").Value = 0.7023186840741513
$ws.Range("Q2").Value = 475.8763710750933
$ws.Range("R2").Value = 4282.88733967584
$ws.Range("S2").Value = 0.3727735200372571
$ws.Range("T2").Value = 0.3727735200372571
$ws.Range("G3").Value = 6.775549333333333
$ws.Range("H3").Value = 20.326648
$ws.Range("I3").Value = 0.5307754563424079
$ws.Range("J3").Value = 0.5307754563424079
$ws.Range("O3").Value = 0.06551129587759326
$ws.Range("P3").Value = 0.06551129587759325
$ws.Range("Q3").Value = 44.38907643152533
$ws.Range("R3").Value = 399.501687883728
$ws.Range("S3").Value = 0.03477178796501206
$ws.Range("T3").Value = 0.03477178796501205
$ws.Range("G4").Value = 6.775549333333333
$ws.Range("H4").Value = 20.326648
$ws.Range("I4").Value = 0.5307754563424079
$ws.Range("J4").Value = 0.5307754563424079
$ws.Range("M4").Value = 4.268944666666666
$ws.Range("N4").Value = 12.806834
$ws.Range("O4").Value = 0.04268793224112385
$ws.Range("P4").Value = 0.04268793224112385
$ws.Range("Q4").Value = 28.92444519027022
$ws.Range("R4").Value = 260.320006712432
$ws.Range("S4").Value = 0.0226577067155963
$ws.Range("T4").Value = 0.0226577067155963
$ws.Range("G5").Value = 6.775549333333333
$ws.Range("H5").Value = 20.326648
$ws.Range("I5").Value = 0.5307754563424079
$ws.Range("J5").Value = 0.5307754563424079
$ws.Range("M5").Value = 18.948881
$ws.Range("N5").Value = 56.846643
$ws.Range("O5").Value = 0.1894820878071316
$ws.Range("P5").Value = 0.1894820878071315
$ws.Range("Q5").Value = 128.3890780269627
$ws.Range("R5").Value = 1155.501702242664
$ws.Range("S5").Value = 0.1005724416245425
$ws.Range("T5").Value = 0.1005724416245424
$ws.Range("I6").Value = 0.3421215311185197
$ws.Range("J6").Value = 0.3421215311185197
$ws.Range("M6").Value = 70.23436
$ws.Range("N6").Value = 210.70308
$ws.Range("O6").Value = 0.7023186840741513
$ws.Range("P6").Value = 0.7023186840741513
$ws.Range("Q6").Value = 306.7352695945066
$ws.Range("R6").Value = 2760.61742635056
$ws.Range("S6").Value = 0.2402783435285926
$ws.Range("T6").Value = 0.2402783435285926
$ws.Range("I7").Value = 0.3421215311185197
$ws.Range("J7").Value = 0.3421215311185197
$ws.Range("O7").Value = 0.06551129587759326
$ws.Range("P7").Value = 0.06551129587759325
$ws.Range("S7").Value = 0.02241282485120058
$ws.Range("T7").Value = 0.02241282485120057
$ws.Range("I8").Value = 0.3421215311185197
$ws.Range("J8").Value = 0.3421215311185197
$ws.Range("M8").Value = 4.268944666666666
$ws.Range("N8").Value = 12.806834
$ws.Range("O8").Value = 0.04268793224112385
$ws.Range("P8").Value = 0.04268793224112385
$ws.Range("Q8").Value = 18.64380757814311
$ws.Range("R8").Value = 167.794268203288
$ws.Range("S8").Value = 0.01460446073861692
$ws.Range("T8").Value = 0.01460446073861692
$ws.Range("I9").Value = 0.3421215311185197
$ws.Range("J9").Value = 0.3421215311185197
$ws.Range("M9").Value = 18.948881
$ws.Range("N9").Value = 56.846643
$ws.Range("O9").Value = 0.1894820878071316
$ws.Range("P9").Value = 0.1894820878071315
$ws.Range("Q9").Value = 82.75565011269734
$ws.Range("R9").Value = 744.800851014276
$ws.Range("S9").Value = 0.06482590200010965
$ws.Range("T9").Value = 0.06482590200010964
$ws.Range("G10").Value = 1.622518
$ws.Range("H10").Value = 4.867554
$ws.Range("I10").Value = 0.1271030125390725
$ws.Range("J10").Value = 0.1271030125390725
$ws.Range("M10").Value = 70.23436
$ws.Range("N10").Value = 210.70308
$ws.Range("O10").Value = 0.7023186840741513
$ws.Range("P10").Value = 0.7023186840741513
$ws.Range("Q10").Value = 113.95651331848
$ws.Range("R10").Value = 1025.60861986632
$ws.Range("S10").Value = 0.08926682050830174
$ws.Range("T10").Value = 0.08926682050830174
$ws.Range("G11").Value = 1.622518
$ws.Range("H11").Value = 4.867554
$ws.Range("I11").Value = 0.1271030125390725
$ws.Range("J11").Value = 0.1271030125390725
$ws.Range("O11").Value = 0.06551129587759326
$ws.Range("P11").Value = 0.06551129587759325
$ws.Range("Q11").Value = 10.629702769516
$ws.Range("R11").Value = 95.667324925644
$ws.Range("S11").Value = 0.008326683061380624
$ws.Range("T11").Value = 0.008326683061380622
$ws.Range("G12").Value = 1.622518
$ws.Range("H12").Value = 4.867554
$ws.Range("I12").Value = 0.1271030125390725
$ws.Range("J12").Value = 0.1271030125390725
$ws.Range("M12").Value = 4.268944666666666
$ws.Range("N12").Value = 12.806834
$ws.Range("O12").Value = 0.04268793224112385
$ws.Range("P12").Value = 0.04268793224112385
$ws.Range("Q12").Value = 6.926439562670667
$ws.Range("R12").Value = 62.33795606403599
$ws.Range("S12").Value = 0.005425764786910642
$ws.Range("T12").Value = 0.005425764786910642
$ws.Range("G13").Value = 1.622518
$ws.Range("H13").Value = 4.867554
$ws.Range("I13").Value = 0.1271030125390725
$ws.Range("J13").Value = 0.1271030125390725
$ws.Range("M13").Value = 18.948881
$ws.Range("N13").Value = 56.846643
$ws.Range("O13").Value = 0.1894820878071316
$ws.Range("P13").Value = 0.1894820878071315
$ws.Range("Q13").Value = 30.744900502358
$ws.Range("R13").Value = 276.704104521222
$ws.Range("S13").Value = 0.02408374418247948
$ws.Range("T13").Value = 0.02408374418247948
